# Auto-update: PANELES GLOBALES MENSUALES
# Adds two new classifier rows (SFE-MAQE, SFE-QUIM) under the "Industrias"
# folder, right after the existing "SFE-MAQ" row, shifting every following
# row down by two, and refreshes the sheet's view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "SFE-MAQ" / "Industrias" currently lives on row 62 - insert the two new
# classifier rows immediately below it.
$ws.Rows("63:64").Insert()

$ws.Range("A63").Value = "SFE-MAQE"
$ws.Range("B63").Value = "Industrias"

$ws.Range("A64").Value = "SFE-QUIM"
$ws.Range("B64").Value = "Industrias"

# The color-scale conditional format on column B applied to B2:B82 needs to
# grow so it keeps covering the full (now longer) list.
$cf = $ws.Range("B2:B82").FormatConditions.Item(1)
$cf.ModifyAppliesToRange($ws.Range("B2:B84"))

# Restore the saved view/selection state.
[void]$ws.Activate()
[void]$ws.Range("D70").Select()
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 1
